$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author renamed the shared string "CONST1" to "CONSTM1".
# Cell A14 held "CONST1"; update it to "CONSTM1".
$ws.Range("A14").Value = "CONSTM1"

# Update the visible selection to match the post-edit state (C15).
$ws.Range("C15").Select()
